$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1543.921
$ws.Range("J17").Value = 1543.921
$ws.Range("L17").Value = 4631.763
$ws.Range("N17").Value = -4967.763

$ws.Range("H28").Value = 1773.409
$ws.Range("I28").Value = 1397.6666
$ws.Range("K28").Value = 1397.6666
$ws.Range("M28").Value = -912.6666

$ws.Range("H40").Value = 17254674
$ws.Range("I40").Value = 14320.353
$ws.Range("K40").Value = 14320.353
$ws.Range("M40").Value = -14145.353

$ws.Range("H64").Value = 21646106
$ws.Range("I64").Value = 7069400
$ws.Range("J64").Value = 38465380
$ws.Range("K64").Value = 7069400
$ws.Range("L64").Value = 38465380
$ws.Range("M64").Value = -7069152
$ws.Range("N64").Value = -38465876

$ws.Range("H67").Value = 21646106
$ws.Range("I67").Value = 7069400
$ws.Range("J67").Value = 38465380
$ws.Range("K67").Value = 7069400
$ws.Range("L67").Value = 38465380
$ws.Range("M67").Value = -7068542
$ws.Range("N67").Value = -38467096

$ws.Range("H86").Value = 1937736.4
$ws.Range("I86").Value = 2945211.2
$ws.Range("J86").Value = 90699.336
$ws.Range("K86").Value = 2945211.2
$ws.Range("L86").Value = 90699.336
$ws.Range("M86").Value = -2944088.2
$ws.Range("N86").Value = -92945.336

$ws.Range("H89").Value = 1937736.4
$ws.Range("I89").Value = 2945211.2
$ws.Range("J89").Value = 90699.336
$ws.Range("K89").Value = 14726056
$ws.Range("L89").Value = 453496.68
$ws.Range("M89").Value = -14720440
$ws.Range("N89").Value = -464728.68

$ws.Range("H112").Value = 1403.0851
$ws.Range("I112").Value = 1053.75
$ws.Range("J112").Value = 1435.5814
$ws.Range("K112").Value = 3161.25
$ws.Range("L112").Value = 4306.7442
$ws.Range("M112").Value = -2053.25
$ws.Range("N112").Value = -6522.7442

$ws.Range("H135").Value = 2753.9333
$ws.Range("I135").Value = 904.381
$ws.Range("J135").Value = 7069.5557
$ws.Range("K135").Value = 8139.429
$ws.Range("L135").Value = 63626.0013
$ws.Range("M135").Value = -5604.429
$ws.Range("N135").Value = -68696.0013

$ws.Range("H137").Value = 17171860
$ws.Range("I137").Value = 2000978
$ws.Range("J137").Value = 22228820
$ws.Range("K137").Value = 6002934
$ws.Range("L137").Value = 66686460
$ws.Range("M137").Value = -6000384
$ws.Range("N137").Value = -66691560

$ws.Range("H138").Value = 2219.16
$ws.Range("I138").Value = 1877.7778
$ws.Range("J138").Value = 2411.1875
$ws.Range("K138").Value = 5633.3334
$ws.Range("L138").Value = 7233.5625
$ws.Range("M138").Value = -493.3334000000004
$ws.Range("N138").Value = -17513.5625

$ws.Range("H141").Value = 3528.652
$ws.Range("I141").Value = 3670.8635
$ws.Range("K141").Value = 11012.5905
$ws.Range("M141").Value = -5832.5905


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10346.594
$ws.Range("I32").Value = 10436.08
$ws.Range("K32").Value = 10436.08
$ws.Range("M32").Value = -10149.08

$ws.Range("H45").Value = 5712.222
$ws.Range("I45").Value = 5504.3335
$ws.Range("K45").Value = 5504.3335
$ws.Range("M45").Value = -5127.3335

$ws.Range("H74").Value = 5953435
$ws.Range("I74").Value = 8621425
$ws.Range("J74").Value = 1764.7693
$ws.Range("K74").Value = 8621425
$ws.Range("L74").Value = 1764.7693
$ws.Range("M74").Value = -8620551
$ws.Range("N74").Value = -3512.7693

$ws.Range("H77").Value = 5953435
$ws.Range("I77").Value = 8621425
$ws.Range("J77").Value = 1764.7693
$ws.Range("K77").Value = 43107125
$ws.Range("L77").Value = 8823.8465
$ws.Range("M77").Value = -43102757
$ws.Range("N77").Value = -17559.8465

$ws.Range("H132").Value = 17499.342
$ws.Range("J132").Value = 6767.84
$ws.Range("L132").Value = 20303.52
$ws.Range("N132").Value = -25363.52


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8339.267
$ws.Range("J20").Value = 8221.166999999999
$ws.Range("L20").Value = 8221.166999999999
$ws.Range("N20").Value = -8715.166999999999

$ws.Range("H105").Value = 2530.9092
$ws.Range("I105").Value = 2209.7144
$ws.Range("J105").Value = 3093
$ws.Range("K105").Value = 2209.7144
$ws.Range("L105").Value = 3093
$ws.Range("M105").Value = -462.7143999999998
$ws.Range("N105").Value = -6587

$ws.Range("H107").Value = 4620.875
$ws.Range("I107").Value = 5497.8335
$ws.Range("J107").Value = 1990
$ws.Range("K107").Value = 5497.8335
$ws.Range("L107").Value = 1990
$ws.Range("M107").Value = -3577.8335
$ws.Range("N107").Value = -5830

$ws.Range("H134").Value = 1333.5938
$ws.Range("I134").Value = 822.6
$ws.Range("J134").Value = 8998.5
$ws.Range("K134").Value = 2467.8
$ws.Range("L134").Value = 26995.5
$ws.Range("M134").Value = 67.19999999999982
$ws.Range("N134").Value = -32065.5


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6552.375
$ws.Range("I31").Value = 1910.6666
$ws.Range("J31").Value = 7623.5386
$ws.Range("K31").Value = 1910.6666
$ws.Range("L31").Value = 7623.5386
$ws.Range("M31").Value = -1615.6666
$ws.Range("N31").Value = -8213.5386

$ws.Range("H34").Value = 6552.375
$ws.Range("I34").Value = 1910.6666
$ws.Range("J34").Value = 7623.5386
$ws.Range("K34").Value = 1910.6666
$ws.Range("L34").Value = 7623.5386
$ws.Range("M34").Value = -1708.6666
$ws.Range("N34").Value = -8027.5386

$ws.Range("H86").Value = 7330.4287
$ws.Range("I86").Value = 6399.4
$ws.Range("K86").Value = 6399.4
$ws.Range("M86").Value = -5276.4

$ws.Range("H89").Value = 7330.4287
$ws.Range("I89").Value = 6399.4
$ws.Range("K89").Value = 31997
$ws.Range("M89").Value = -26381

$ws.Range("H109").Value = 74998.75
$ws.Range("J109").Value = 74998.75
$ws.Range("L109").Value = 74998.75
$ws.Range("N109").Value = -77078.75


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1996.6666
$ws.Range("J68").Value = 1996.6666
$ws.Range("L68").Value = 5989.9998
$ws.Range("N68").Value = -7611.9998

$ws.Range("H71").Value = 1996.6666
$ws.Range("J71").Value = 1996.6666
$ws.Range("L71").Value = 17969.9994
$ws.Range("N71").Value = -26081.9994


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4137901.2
$ws.Range("I70").Value = 4550791.5
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 4550791.5
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4550521.5
$ws.Range("N70").Value = -9540

$ws.Range("H73").Value = 4137901.2
$ws.Range("I73").Value = 4550791.5
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 4550791.5
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -4549855.5
$ws.Range("N73").Value = -10872

$ws.Range("H113").Value = 1880
$ws.Range("I113").Value = 1850
$ws.Range("K113").Value = 1850
$ws.Range("M113").Value = 320

$ws.Range("H122").Value = 792595.4
$ws.Range("I122").Value = 2753224.8
$ws.Range("J122").Value = 8343.6
$ws.Range("K122").Value = 8259674.399999999
$ws.Range("L122").Value = 25030.8
$ws.Range("M122").Value = -8257224.399999999
$ws.Range("N122").Value = -29930.8

$ws.Range("H132").Value = 107959.4
$ws.Range("I132").Value = 150128.08
$ws.Range("K132").Value = 450384.24
$ws.Range("M132").Value = -447854.24


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7004.7
$ws.Range("I46").Value = 5401.25
$ws.Range("K46").Value = 5401.25
$ws.Range("M46").Value = -5213.25

$ws.Range("H82").Value = 2841765.8
$ws.Range("I82").Value = 5208962.5
$ws.Range("J82").Value = 1129.6
$ws.Range("K82").Value = 5208962.5
$ws.Range("L82").Value = 1129.6
$ws.Range("M82").Value = -5208601.5
$ws.Range("N82").Value = -1851.6

$ws.Range("H85").Value = 2841765.8
$ws.Range("I85").Value = 5208962.5
$ws.Range("J85").Value = 1129.6
$ws.Range("K85").Value = 5208962.5
$ws.Range("L85").Value = 1129.6
$ws.Range("M85").Value = -5207714.5
$ws.Range("N85").Value = -3625.6

$ws.Range("H122").Value = 43483708
$ws.Range("I122").Value = 71432940
$ws.Range("J122").Value = 7122.5557
$ws.Range("K122").Value = 214298820
$ws.Range("L122").Value = 21367.6671
$ws.Range("M122").Value = -214296370
$ws.Range("N122").Value = -26267.6671

$ws.Range("H132").Value = 4920.636
$ws.Range("I132").Value = 4579.375
$ws.Range("K132").Value = 13738.125
$ws.Range("M132").Value = -11208.125


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3180000.2
$ws.Range("I5").Value = 10000001
$ws.Range("J5").Value = 1475000
$ws.Range("K5").Value = 10000001
$ws.Range("L5").Value = 1475000
$ws.Range("M5").Value = -9999889
$ws.Range("N5").Value = -1475224

$ws.Range("H113").Value = 963.2
$ws.Range("I113").Value = 967.9
$ws.Range("J113").Value = 953.8
$ws.Range("K113").Value = 2903.7
$ws.Range("L113").Value = 2861.4
$ws.Range("M113").Value = -733.6999999999998
$ws.Range("N113").Value = -7201.4

$ws.Range("H116").Value = 108999
$ws.Range("J116").Value = 108999
$ws.Range("L116").Value = 108999
$ws.Range("N116").Value = -118177

$ws.Range("H126").Value = 4417.7856
$ws.Range("I126").Value = 4245
$ws.Range("K126").Value = 12735
$ws.Range("M126").Value = -10265

$ws.Range("H132").Value = 9437014
$ws.Range("I132").Value = 1043.303
$ws.Range("K132").Value = 3129.909000000001
$ws.Range("M132").Value = -599.9090000000006

